$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3852, 4145, 4145, 4400, 4660, 4660, 4660, 4937, 4937, 4966, 4966, 4966, 5032, 5093)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
